$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update "max_dias" (column N) for rows 3-70: most become 7; row 7 becomes the exception (1) ---
$ws.Cells.Item(3, 14).Value = 7
$ws.Cells.Item(4, 14).Value = 7
$ws.Cells.Item(5, 14).Value = 7
$ws.Cells.Item(6, 14).Value = 7
$ws.Cells.Item(7, 14).Value = 1
$ws.Cells.Item(8, 14).Value = 7
$ws.Cells.Item(9, 14).Value = 7
$ws.Cells.Item(10, 14).Value = 7
$ws.Cells.Item(11, 14).Value = 7
$ws.Cells.Item(12, 14).Value = 7
$ws.Cells.Item(13, 14).Value = 7
$ws.Cells.Item(14, 14).Value = 7
$ws.Cells.Item(15, 14).Value = 7
$ws.Cells.Item(16, 14).Value = 7
$ws.Cells.Item(17, 14).Value = 7
$ws.Cells.Item(19, 14).Value = 7
$ws.Cells.Item(22, 14).Value = 7
$ws.Cells.Item(23, 14).Value = 7
$ws.Cells.Item(24, 14).Value = 7
$ws.Cells.Item(25, 14).Value = 7
$ws.Cells.Item(26, 14).Value = 7
$ws.Cells.Item(27, 14).Value = 7
$ws.Cells.Item(28, 14).Value = 7
$ws.Cells.Item(30, 14).Value = 7
$ws.Cells.Item(31, 14).Value = 7
$ws.Cells.Item(33, 14).Value = 7
$ws.Cells.Item(34, 14).Value = 7
$ws.Cells.Item(35, 14).Value = 7
$ws.Cells.Item(36, 14).Value = 7
$ws.Cells.Item(37, 14).Value = 7
$ws.Cells.Item(38, 14).Value = 7
$ws.Cells.Item(39, 14).Value = 7
$ws.Cells.Item(40, 14).Value = 7
$ws.Cells.Item(43, 14).Value = 7
$ws.Cells.Item(45, 14).Value = 7
$ws.Cells.Item(46, 14).Value = 7
$ws.Cells.Item(47, 14).Value = 7
$ws.Cells.Item(48, 14).Value = 7
$ws.Cells.Item(49, 14).Value = 7
$ws.Cells.Item(50, 14).Value = 7
$ws.Cells.Item(51, 14).Value = 7
$ws.Cells.Item(52, 14).Value = 7
$ws.Cells.Item(53, 14).Value = 7
$ws.Cells.Item(54, 14).Value = 7
$ws.Cells.Item(55, 14).Value = 7
$ws.Cells.Item(59, 14).Value = 7
$ws.Cells.Item(60, 14).Value = 7
$ws.Cells.Item(61, 14).Value = 7
$ws.Cells.Item(62, 14).Value = 7
$ws.Cells.Item(63, 14).Value = 7
$ws.Cells.Item(64, 14).Value = 7
$ws.Cells.Item(65, 14).Value = 7
$ws.Cells.Item(66, 14).Value = 7
$ws.Cells.Item(67, 14).Value = 7
$ws.Cells.Item(68, 14).Value = 7
$ws.Cells.Item(69, 14).Value = 7
$ws.Cells.Item(70, 14).Value = 7

# --- Rows 71-86: ids renumbered/shifted, two foods spliced in (amendoim, carne) before ovo,
# and three foods appended after frango (atum, azeite, manteiga) ---
# row 71: banana
$ws.Cells.Item(71, 1).Value = 175
$ws.Cells.Item(71, 2).Value = 'banana'
$ws.Cells.Item(71, 3).Value = 128
$ws.Cells.Item(71, 4).Value = 33.7
$ws.Cells.Item(71, 5).Value = 1.4
$ws.Cells.Item(71, 6).Value = 0.2
$ws.Cells.Item(71, 7).Value = 0.3
$ws.Cells.Item(71, 8).Value = 24
$ws.Cells.Item(71, 9).Value = 15.7
$ws.Cells.Item(71, 10).Value = 0.2
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 1.65
$ws.Cells.Item(71, 13).Value = 10
$ws.Cells.Item(71, 14).Value = 7
$ws.Cells.Item(71, 15).Value = 2

# row 72: achocolatado
$ws.Cells.Item(72, 1).Value = 491
$ws.Cells.Item(72, 2).Value = 'achocolatado'
$ws.Cells.Item(72, 3).Value = 401
$ws.Cells.Item(72, 4).Value = 91.2
$ws.Cells.Item(72, 5).Value = 4.2
$ws.Cells.Item(72, 6).Value = 2.2
$ws.Cells.Item(72, 7).Value = 5.4
$ws.Cells.Item(72, 8).Value = 77
$ws.Cells.Item(72, 9).Value = 0
$ws.Cells.Item(72, 10).Value = 1
$ws.Cells.Item(72, 11).Value = 65
$ws.Cells.Item(72, 12).Value = 2.67
$ws.Cells.Item(72, 13).Value = 4
$ws.Cells.Item(72, 14).Value = 7
$ws.Cells.Item(72, 15).Value = 1

# row 73: acucarCristal
$ws.Cells.Item(73, 1).Value = 492
$ws.Cells.Item(73, 2).Value = 'acucarCristal'
$ws.Cells.Item(73, 3).Value = 387
$ws.Cells.Item(73, 4).Value = 99.6
$ws.Cells.Item(73, 5).Value = 0.3
$ws.Cells.Item(73, 6).Value = 0
$ws.Cells.Item(73, 7).Value = 0.2
$ws.Cells.Item(73, 8).Value = 1
$ws.Cells.Item(73, 9).Value = 0
$ws.Cells.Item(73, 10).Value = 0
$ws.Cells.Item(73, 11).Value = 0
$ws.Cells.Item(73, 12).Value = 0.5
$ws.Cells.Item(73, 13).Value = 3
$ws.Cells.Item(73, 14).Value = 7
$ws.Cells.Item(73, 15).Value = 1

# row 74: acucarMascavo
$ws.Cells.Item(74, 1).Value = 493
$ws.Cells.Item(74, 2).Value = 'acucarMascavo'
$ws.Cells.Item(74, 3).Value = 369
$ws.Cells.Item(74, 4).Value = 94.5
$ws.Cells.Item(74, 5).Value = 0.8
$ws.Cells.Item(74, 6).Value = 0.1
$ws.Cells.Item(74, 7).Value = 8.3
$ws.Cells.Item(74, 8).Value = 80
$ws.Cells.Item(74, 9).Value = 0
$ws.Cells.Item(74, 10).Value = 0.5
$ws.Cells.Item(74, 11).Value = 25
$ws.Cells.Item(74, 12).Value = 2.179
$ws.Cells.Item(74, 13).Value = 3
$ws.Cells.Item(74, 14).Value = 7
$ws.Cells.Item(74, 15).Value = 1

# row 75: acucarRefinado
$ws.Cells.Item(75, 1).Value = 494
$ws.Cells.Item(75, 2).Value = 'acucarRefinado'
$ws.Cells.Item(75, 3).Value = 387
$ws.Cells.Item(75, 4).Value = 99.5
$ws.Cells.Item(75, 5).Value = 0.3
$ws.Cells.Item(75, 6).Value = 0
$ws.Cells.Item(75, 7).Value = 0.1
$ws.Cells.Item(75, 8).Value = 1
$ws.Cells.Item(75, 9).Value = 0
$ws.Cells.Item(75, 10).Value = 0
$ws.Cells.Item(75, 11).Value = 12
$ws.Cells.Item(75, 12).Value = 0.49
$ws.Cells.Item(75, 13).Value = 3
$ws.Cells.Item(75, 14).Value = 7
$ws.Cells.Item(75, 15).Value = 1

# row 76: chocolateAoLeite
$ws.Cells.Item(76, 1).Value = 495
$ws.Cells.Item(76, 2).Value = 'chocolateAoLeite'
$ws.Cells.Item(76, 3).Value = 540
$ws.Cells.Item(76, 4).Value = 59.6
$ws.Cells.Item(76, 5).Value = 7.2
$ws.Cells.Item(76, 6).Value = 30.3
$ws.Cells.Item(76, 7).Value = 1.6
$ws.Cells.Item(76, 8).Value = 57
$ws.Cells.Item(76, 9).Value = 0
$ws.Cells.Item(76, 10).Value = 1.1
$ws.Cells.Item(76, 11).Value = 77
$ws.Cells.Item(76, 12).Value = 5.99
$ws.Cells.Item(76, 13).Value = 4
$ws.Cells.Item(76, 14).Value = 7
$ws.Cells.Item(76, 15).Value = 1

# row 77: chocolateMeioAmargo
$ws.Cells.Item(77, 1).Value = 498
$ws.Cells.Item(77, 2).Value = 'chocolateMeioAmargo'
$ws.Cells.Item(77, 3).Value = 475
$ws.Cells.Item(77, 4).Value = 62.4
$ws.Cells.Item(77, 5).Value = 4.9
$ws.Cells.Item(77, 6).Value = 29.9
$ws.Cells.Item(77, 7).Value = 3.6
$ws.Cells.Item(77, 8).Value = 107
$ws.Cells.Item(77, 9).Value = 2.1
$ws.Cells.Item(77, 10).Value = 1.5
$ws.Cells.Item(77, 11).Value = 9
$ws.Cells.Item(77, 12).Value = 11.29
$ws.Cells.Item(77, 13).Value = 4
$ws.Cells.Item(77, 14).Value = 7
$ws.Cells.Item(77, 15).Value = 1

# row 78: cocada
$ws.Cells.Item(78, 1).Value = 499
$ws.Cells.Item(78, 2).Value = 'cocada'
$ws.Cells.Item(78, 3).Value = 449
$ws.Cells.Item(78, 4).Value = 81.4
$ws.Cells.Item(78, 5).Value = 1.1
$ws.Cells.Item(78, 6).Value = 13.6
$ws.Cells.Item(78, 7).Value = 1.2
$ws.Cells.Item(78, 8).Value = 17
$ws.Cells.Item(78, 9).Value = 0
$ws.Cells.Item(78, 10).Value = 0.4
$ws.Cells.Item(78, 11).Value = 29
$ws.Cells.Item(78, 12).Value = 8.39
$ws.Cells.Item(78, 13).Value = 3
$ws.Cells.Item(78, 14).Value = 7
$ws.Cells.Item(78, 15).Value = 1

# row 79: feijao
$ws.Cells.Item(79, 1).Value = 561
$ws.Cells.Item(79, 2).Value = 'feijao'
$ws.Cells.Item(79, 3).Value = 76
$ws.Cells.Item(79, 4).Value = 13.6
$ws.Cells.Item(79, 5).Value = 4.8
$ws.Cells.Item(79, 6).Value = 0.5
$ws.Cells.Item(79, 7).Value = 1.3
$ws.Cells.Item(79, 8).Value = 42
$ws.Cells.Item(79, 9).Value = 0
$ws.Cells.Item(79, 10).Value = 0.7
$ws.Cells.Item(79, 11).Value = 2
$ws.Cells.Item(79, 12).Value = 0.699
$ws.Cells.Item(79, 13).Value = 20
$ws.Cells.Item(79, 14).Value = 7
$ws.Cells.Item(79, 15).Value = 6

# row 80: amendoim
$ws.Cells.Item(80, 1).Value = 558
$ws.Cells.Item(80, 2).Value = 'amendoim'
$ws.Cells.Item(80, 3).Value = 606
$ws.Cells.Item(80, 4).Value = 18.7
$ws.Cells.Item(80, 5).Value = 22.5
$ws.Cells.Item(80, 6).Value = 54
$ws.Cells.Item(80, 7).Value = 1.3
$ws.Cells.Item(80, 8).Value = 159
$ws.Cells.Item(80, 9).Value = 0
$ws.Cells.Item(80, 10).Value = 2.1
$ws.Cells.Item(80, 11).Value = 376
$ws.Cells.Item(80, 12).Value = 3.89
$ws.Cells.Item(80, 13).Value = 2
$ws.Cells.Item(80, 14).Value = 7
$ws.Cells.Item(80, 15).Value = 6

# row 81: carne
$ws.Cells.Item(81, 1).Value = 377
$ws.Cells.Item(81, 2).Value = 'carne'
$ws.Cells.Item(81, 3).Value = 219
$ws.Cells.Item(81, 4).Value = 0
$ws.Cells.Item(81, 5).Value = 35.9
$ws.Cells.Item(81, 6).Value = 7.3
$ws.Cells.Item(81, 7).Value = 3
$ws.Cells.Item(81, 8).Value = 27
$ws.Cells.Item(81, 9).Value = 0
$ws.Cells.Item(81, 10).Value = 8.1
$ws.Cells.Item(81, 11).Value = 60
$ws.Cells.Item(81, 12).Value = 5.69
$ws.Cells.Item(81, 13).Value = 20
$ws.Cells.Item(81, 14).Value = 7
$ws.Cells.Item(81, 15).Value = 5

# row 82: ovo
$ws.Cells.Item(82, 1).Value = 488
$ws.Cells.Item(82, 2).Value = 'ovo'
$ws.Cells.Item(82, 3).Value = 146
$ws.Cells.Item(82, 4).Value = 0.6
$ws.Cells.Item(82, 5).Value = 13.3
$ws.Cells.Item(82, 6).Value = 9.5
$ws.Cells.Item(82, 7).Value = 1.5
$ws.Cells.Item(82, 8).Value = 11
$ws.Cells.Item(82, 9).Value = 0
$ws.Cells.Item(82, 10).Value = 1.2
$ws.Cells.Item(82, 11).Value = 146
$ws.Cells.Item(82, 12).Value = 1.8
$ws.Cells.Item(82, 13).Value = 15
$ws.Cells.Item(82, 14).Value = 7
$ws.Cells.Item(82, 15).Value = 5

# row 83: frango
$ws.Cells.Item(83, 1).Value = 406
$ws.Cells.Item(83, 2).Value = 'frango'
$ws.Cells.Item(83, 3).Value = 212
$ws.Cells.Item(83, 4).Value = 0
$ws.Cells.Item(83, 5).Value = 33.4
$ws.Cells.Item(83, 6).Value = 7.6
$ws.Cells.Item(83, 7).Value = 0.5
$ws.Cells.Item(83, 8).Value = 18
$ws.Cells.Item(83, 9).Value = 0
$ws.Cells.Item(83, 10).Value = 1
$ws.Cells.Item(83, 11).Value = 56
$ws.Cells.Item(83, 12).Value = 2.49
$ws.Cells.Item(83, 13).Value = 20
$ws.Cells.Item(83, 14).Value = 7
$ws.Cells.Item(83, 15).Value = 5

# row 84: atum
$ws.Cells.Item(84, 1).Value = 277
$ws.Cells.Item(84, 2).Value = 'atum'
$ws.Cells.Item(84, 3).Value = 166
$ws.Cells.Item(84, 4).Value = 0
$ws.Cells.Item(84, 5).Value = 26.2
$ws.Cells.Item(84, 6).Value = 6
$ws.Cells.Item(84, 7).Value = 1.2
$ws.Cells.Item(84, 8).Value = 29
$ws.Cells.Item(84, 9).Value = 0
$ws.Cells.Item(84, 10).Value = 0.6
$ws.Cells.Item(84, 11).Value = 362
$ws.Cells.Item(84, 12).Value = 5.35
$ws.Cells.Item(84, 13).Value = 2
$ws.Cells.Item(84, 14).Value = 7
$ws.Cells.Item(84, 15).Value = 5

# row 85: azeite
$ws.Cells.Item(85, 1).Value = 260
$ws.Cells.Item(85, 2).Value = 'azeite'
$ws.Cells.Item(85, 3).Value = 884
$ws.Cells.Item(85, 4).Value = 0
$ws.Cells.Item(85, 5).Value = 0
$ws.Cells.Item(85, 6).Value = 100
$ws.Cells.Item(85, 7).Value = 0
$ws.Cells.Item(85, 8).Value = 0
$ws.Cells.Item(85, 9).Value = 0
$ws.Cells.Item(85, 10).Value = 0
$ws.Cells.Item(85, 11).Value = 0
$ws.Cells.Item(85, 12).Value = 11.19
$ws.Cells.Item(85, 13).Value = 1
$ws.Cells.Item(85, 14).Value = 7
$ws.Cells.Item(85, 15).Value = 7

# row 86: manteiga
$ws.Cells.Item(86, 1).Value = 261
$ws.Cells.Item(86, 2).Value = 'manteiga'
$ws.Cells.Item(86, 3).Value = 726
$ws.Cells.Item(86, 4).Value = 0.1
$ws.Cells.Item(86, 5).Value = 0.4
$ws.Cells.Item(86, 6).Value = 82.4
$ws.Cells.Item(86, 7).Value = 0.2
$ws.Cells.Item(86, 8).Value = 1
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 10).Value = 0
$ws.Cells.Item(86, 11).Value = 579
$ws.Cells.Item(86, 12).Value = 6.49
$ws.Cells.Item(86, 13).Value = 1
$ws.Cells.Item(86, 14).Value = 7
$ws.Cells.Item(86, 15).Value = 7

# --- Freeze header row (row 1) and restore selection/view state ---
$ws.Range("A1").Select()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("Q7").Select()
